$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new address row (row 3)
$ws.Range("A3").Value = "26772 Calle Maria"
$ws.Range("D3").Value = "Capistrano Beach"
$ws.Range("E3").Value = "ca"
$ws.Range("F3").Value = "us"
$ws.Range("G3").Value = 92624

# Update the selection to match the committed state
$ws.Range("H3").Select()
